$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Errori")

# Header + the two "commesse tassative" rows flagged as problematic
$ws.Range("A1").Value = "id"
$ws.Range("A2").Value = 252681
$ws.Range("A3").Value = 253497

# Widen column A so the ids are readable (stored width="30" in the xlsx)
$ws.Columns.Item(1).ColumnWidth = 29.1666666666667
